$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Top 50 Cryptocurrencies" ---
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$r2 = @(2,'Bitcoin','btc',98953,1960055536355,117216670870,2.03171)
$r3 = @(3,'Ethereum','eth',3405.98,410371174804,56012115569,8.61162)
$r4 = @(4,'Tether','usdt',1.001,130826999402,115551101401,-0.16443)
$r5 = @(5,'Solana','sol',261.32,124152554628,15095440759,8.61017)
$r6 = @(6,'BNB','bnb',636.05,92809270596,2499467256,4.31081)
$r7 = @(7,'XRP','xrp',1.4,79939690866,17839941365,25.65458)
$r8 = @(8,'Dogecoin','doge',0.397586,58438525648,10229008017,2.55081)
$r9 = @(9,'USDC','usdc',0.999424,38293352460,16171222892,-0.19814)
$r10 = @(10,'Lido Staked Ether','steth',3404.14,33379470053,148043344,8.945650000000001)
$r11 = @(11,'Cardano','ada',0.89227,31956726304,3233777838,11.80011)
$r12 = @(12,'TRON','trx',0.200695,17339715807,1096809323,1.56179)
$r13 = @(13,'Avalanche','avax',36.52,14954843621,1047754060,6.78292)
$r14 = @(14,'Shiba Inu','shib',0.00002515,14807481931,1617812388,4.32539)
$r15 = @(15,'Wrapped stETH','wsteth',4005.5,14464595442,167440045,7.863)
$r16 = @(16,'Wrapped Bitcoin','wbtc',98855,14370685455,855965571,2.29304)
$r17 = @(17,'Toncoin','ton',5.59,14225472670,631030859,4.16661)
$r18 = @(18,'Sui','sui',3.64,10349961705,2382450464,1.06058)
$r19 = @(19,'Bitcoin Cash','bch',497.01,9844997186,2191495355,1.75481)
$r20 = @(20,'WETH','weth',3405.53,9710866412,2043740645,8.840859999999999)
$r21 = @(21,'Chainlink','link',15.34,9630736423,1247912196,5.66124)
$r22 = @(22,'Polkadot','dot',6.25,9011225759,820684345,9.504149999999999)
$r23 = @(23,'Pepe','pepe',0.0000214,8995247135,7004707459,9.67314)
$r24 = @(24,'Stellar','xlm',0.286757,8635291214,2342672074,19.95186)
$r25 = @(25,'LEO Token','leo',8.779999999999999,8092311033,3449952,2.85493)
$r26 = @(26,'NEAR Protocol','near',5.84,7113376462,1011267993,4.81004)
$r27 = @(27,'Litecoin','ltc',91,6853292888,1455521572,5.54842)
$r28 = @(28,'Aptos','apt',12.2,6501949665,887636727,3.75529)
$r29 = @(29,'Wrapped eETH','weeth',3595.01,6246701016,103196739,9.221399999999999)
$r30 = @(30,'Uniswap','uni',9.48,5684456629,859274420,7.63452)
$r31 = @(31,'Cronos','cro',0.195162,5275716632,118472792,10.57972)
$r32 = @(32,'USDS','usds',0.997936,5222818939,16029933,-0.41101)
$r33 = @(33,'Hedera','hbar',0.133256,5116781830,887310034,6.06212)
$r34 = @(34,'Internet Computer','icp',9.73,4617180476,271410334,7.07758)
$r35 = @(35,'Ethereum Classic','etc',28.15,4216502883,897048293,6.7889)
$r36 = @(36,'Bonk','bonk',0.00005238,3933453185,1758900842,2.50637)
$r37 = @(37,'Render','render',7.45,3856984122,442272348,0.20977)
$r38 = @(38,'Kaspa','kas',0.151111,3808008207,153886307,0.23278)
$r39 = @(39,'POL (ex-MATIC)','pol',0.475975,3800453937,486972368,7.86641)
$r40 = @(40,'Bittensor','tao',509.38,3759806505,258992108,3.4544)
$r41 = @(41,'Ethena USDe','usde',1.001,3686526730,227406327,-0.19774)
$r42 = @(42,'WhiteBIT Coin','wbt',24.8,3573124265,38721916,2.67548)
$r43 = @(43,'Dai','dai',0.999326,3441593689,160724019,-0.22204)
$r44 = @(44,'MANTRA','om',3.8,3435649852,305465958,4.68595)
$r45 = @(45,'dogwifhat','wif',3.4,3397265421,1279552228,5.6871)
$r46 = @(46,'Artificial Superintelligence Alliance','fet',1.29,3360582858,476746349,3.77132)
$r47 = @(47,'Arbitrum','arb',0.793611,3252551349,1697016111,12.86694)
$r48 = @(48,'Monero','xmr',160.53,2961967862,84034398,-1.04901)
$r49 = @(49,'Stacks','stx',1.96,2938800555,385861546,1.97972)
$r50 = @(50,'Filecoin','fil',4.72,2837726335,592441880,8.15019)
$r51 = @(51,'OKB','okb',46.73,2806970427,20013929,5.40261)

$cryptoData = @($r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12, $r13, $r14, $r15, $r16, $r17, $r18, $r19, $r20, $r21, $r22, $r23, $r24, $r25, $r26, $r27, $r28, $r29, $r30, $r31, $r32, $r33, $r34, $r35, $r36, $r37, $r38, $r39, $r40, $r41, $r42, $r43, $r44, $r45, $r46, $r47, $r48, $r49, $r50, $r51)

foreach ($row in $cryptoData) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
}

# --- Sheet 2: "Top 5 by Market Cap" ---
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws2.Range("B2").Value = 1960055536355
$ws2.Range("B3").Value = 410371174804
$ws2.Range("B4").Value = 130826999402
$ws2.Range("B5").Value = 124152554628
$ws2.Range("B6").Value = 92809270596

# --- Sheet 3: "Summary" ---
$ws3 = $wb.Worksheets.Item("Summary")
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = '$4360.50'
$ws3.Range("B3").Value = "XRP (25.65%)"
$ws3.Range("B4").Value = "Monero (-1.05%)"
